$wb = $excel.ActiveWorkbook

# --- Sheet "FEINmismatch": rows 2-30 ---
$ws1 = $wb.Worksheets.Item("FEINmismatch")

# Column A (Result) - all rows become "Pass"
$ws1.Range("A2:A30").Value = "Pass"

# Column B (Date) - update Execute timestamps to the new test run values
$sheet1Dates = @(
    "Mon Apr 08 17:24:32 EDT 2024",
    "Mon Apr 08 17:24:46 EDT 2024",
    "Mon Apr 08 17:24:57 EDT 2024",
    "Mon Apr 08 17:25:08 EDT 2024",
    "Mon Apr 08 17:25:19 EDT 2024",
    "Mon Apr 08 17:25:30 EDT 2024",
    "Mon Apr 08 17:25:41 EDT 2024",
    "Mon Apr 08 17:25:52 EDT 2024",
    "Mon Apr 08 17:26:04 EDT 2024",
    "Mon Apr 08 17:26:16 EDT 2024",
    "Mon Apr 08 17:26:27 EDT 2024",
    "Mon Apr 08 17:26:37 EDT 2024",
    "Mon Apr 08 17:26:49 EDT 2024",
    "Mon Apr 08 17:27:00 EDT 2024",
    "Mon Apr 08 17:27:11 EDT 2024",
    "Mon Apr 08 17:27:22 EDT 2024",
    "Mon Apr 08 17:27:33 EDT 2024",
    "Mon Apr 08 17:27:44 EDT 2024",
    "Mon Apr 08 17:27:55 EDT 2024",
    "Mon Apr 08 17:28:05 EDT 2024",
    "Mon Apr 08 17:28:16 EDT 2024",
    "Mon Apr 08 17:28:27 EDT 2024",
    "Mon Apr 08 17:28:38 EDT 2024",
    "Mon Apr 08 17:28:49 EDT 2024",
    "Mon Apr 08 17:29:00 EDT 2024",
    "Mon Apr 08 17:29:10 EDT 2024",
    "Mon Apr 08 17:29:21 EDT 2024",
    "Mon Apr 08 17:29:32 EDT 2024",
    "Mon Apr 08 17:29:43 EDT 2024"
)

for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

# --- Sheet "FEINSSNmismatch": rows 2-19 ---
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

$sheet2Dates = @(
    "Mon Apr 08 17:29:54 EDT 2024",
    "Mon Apr 08 17:30:05 EDT 2024",
    "Mon Apr 08 17:30:15 EDT 2024",
    "Mon Apr 08 17:30:26 EDT 2024",
    "Mon Apr 08 17:30:36 EDT 2024",
    "Mon Apr 08 17:30:47 EDT 2024",
    "Mon Apr 08 17:30:57 EDT 2024",
    "Mon Apr 08 17:31:08 EDT 2024",
    "Mon Apr 08 17:31:19 EDT 2024",
    "Mon Apr 08 17:31:29 EDT 2024",
    "Mon Apr 08 17:31:39 EDT 2024",
    "Mon Apr 08 17:31:50 EDT 2024",
    "Mon Apr 08 17:32:00 EDT 2024",
    "Mon Apr 08 17:32:11 EDT 2024",
    "Mon Apr 08 17:32:21 EDT 2024",
    "Mon Apr 08 17:32:32 EDT 2024",
    "Mon Apr 08 17:32:42 EDT 2024",
    "Mon Apr 08 17:32:53 EDT 2024"
)

for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}

Write-Host "Updated RAD EL-Motor Fuel Tax results."
